# Fix "Total" column formatting on the Students sheet (append ".0" to the
# integer part of each total, e.g. "22 (71.0%)" -> "22.0 (71.0%)") and
# correct the stale Points value on the Points sheet (A2: 3 -> 4) so it
# matches the recomputed "Total" for Student E (4 (12.9%)).

$wb = $excel.ActiveWorkbook

# --- Students sheet: column D ("Total"), rows 2-11 ---
$studentsWs = $wb.Worksheets.Item("Students")

$totals = @{
    2  = "22.0 (71.0%)"
    3  = "27.0 (87.1%)"
    4  = "17.0 (54.8%)"
    5  = "9.0 (29.0%)"
    6  = "4.0 (12.9%)"
    7  = "27.0 (87.1%)"
    8  = "13.0 (41.9%)"
    9  = "19.0 (61.3%)"
    10 = "20.0 (64.5%)"
    11 = "22.0 (71.0%)"
}

foreach ($row in $totals.Keys) {
    $studentsWs.Range("D$row").Value = $totals[$row]
}

# --- Points sheet: A2, 3 -> 4 ---
$pointsWs = $wb.Worksheets.Item("Points")
$pointsWs.Range("A2").Value = 4
